$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J20").Value = "G418"
$ws.Range("J21").Value = "G418"
$ws.Range("J22").Value = "G418"

$ws.Range("K22").Select()
